$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.459.58'
$ws.Range("E2").Value = '  +0.51%  '

$ws.Range("D3").Value = '2.688.43'
$ws.Range("E3").Value = '  +1.74%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.05'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.71'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.63%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D9").Value = '2.687.36'
$ws.Range("E9").Value = '  +1.73%  '

$ws.Range("E10").Value = '  -0.40%  '

$ws.Range("E11").Value = '  -0.66%  '

$ws.Range("E12").Value = '  +1.10%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.359'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.42%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.16'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.81%  '

$ws.Range("D15").Value = '3.177.82'
$ws.Range("E15").Value = '  +1.71%  '

$ws.Range("E16").Value = '  -0.35%  '

$ws.Range("D17").Value = '68.403.97'
$ws.Range("E17").Value = '  +0.38%  '

$ws.Range("D18").Value = '2.689.25'
$ws.Range("E18").Value = '  +2.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.84'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.44%  '

$ws.Range("E20").Value = '  +1.04%  '

$ws.Range("E21").Value = '  +3.28%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.53'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.98%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.87'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.97%  '

$ws.Range("E24").Value = '  +2.57%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.46'
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.99'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.73%  '

$ws.Range("D28").Value = '2.824.91'
$ws.Range("E28").Value = '  +1.84%  '

$ws.Range("E29").Value = '  +0.34%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.12%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '570.21'
$ws.Range("D31").Style = "Normal"

$ws.Range("E32").Value = '  +2.72%  '

$ws.Range("E33").Value = '  +3.81%  '

$ws.Range("E34").Value = '  +5.31%  '

$ws.Range("E36").Value = '  +6.49%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.05%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '161.85'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.45%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.83'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.01%  '

$ws.Range("E40").Value = '  +1.84%  '

$ws.Range("E41").Value = '  +2.08%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.39'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.46%  '

$ws.Range("E44").Value = '  +1.69%  '

$ws.Range("E45").Value = '  +0.02%  '

$ws.Range("E46").Value = '  -6.65%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '157.17'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.36%  '

$ws.Range("E48").Value = '  +7.38%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.76'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.55%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.596'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.18%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.97'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.03%  '
